# Add a new "Netherlands" market sheet to the workbook, modeled on the
# existing "Spain" sheet (the last sheet in the gallery), then fill in
# its market-specific data.

$wb = $excel.ActiveWorkbook

# The "Spain" sheet is the template for every per-market sheet in this
# workbook (same layout/styles), so duplicate it to create the new sheet.
$source = $wb.Worksheets.Item("Spain")

# Mirror the "select everything" state the source tab ends up with once it
# stops being the active tab.
$source.Select() | Out-Null
$source.Range("A1:XFD1048576").Select() | Out-Null

# Copy placed immediately after the source sheet -> lands at the end of
# the sheet list.
$source.Copy($null, $source)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Netherlands"

# Fill in the Netherlands-specific values (same cells Germany/Belgium/../
# Spain all use for their market name + user story reference).
$newSheet.Range("B2").Value = "Netherlands Market"
$newSheet.Range("B4").Value = "NGC-3144/T2175"

# Leave the new sheet active/selected, with B4 as the active cell - matches
# how the new sheet was left selected after data entry.
$newSheet.Select() | Out-Null
$newSheet.Range("B4").Select() | Out-Null
